# =====================================================================
# Scheduled market-data refresh for Kujata_Profits workbook
#
# This script pushes freshly-pulled Universalis market prices into the
# per-class profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Only the price / profit columns (H-N) that changed since the last run
# are touched; leve metadata columns (A-G) are left untouched.
#
# A handful of rows gain or lose an HQ-profit figure (column N, and in
# one case column M) when the HQ recipe stops/starts being craftable at
# a profit this cycle - those cells are cleared (or written for the
# first time) with ClearContents()/Value so the cell element is removed
# from / added to the sheet exactly like the source update does.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------- Sheet: ALC ----------------
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 10109.8
$ws.Range("I6").Value = 10109.8
$ws.Range("K6").Value = 30329.4
$ws.Range("M6").Value = -30217.4
# Row 18
$ws.Range("H18").Value = 519.3077
$ws.Range("I18").Value = 250.09091
$ws.Range("K18").Value = 250.09091
$ws.Range("M18").Value = 33.90908999999999
# Row 19
$ws.Range("H19").Value = 415.62964
$ws.Range("J19").Value = 355.0625
$ws.Range("L19").Value = 355.0625
$ws.Range("N19").Value = -705.0625
# Row 116
$ws.Range("H116").Value = 2556.6875
$ws.Range("I116").Value = 2175.5
$ws.Range("J116").Value = 2937.875
$ws.Range("K116").Value = 2175.5
$ws.Range("L116").Value = 2937.875
$ws.Range("M116").Value = 1266.5
$ws.Range("N116").Value = -9821.875
# Row 135
$ws.Range("H135").Value = 27778558
$ws.Range("I135").Value = 493.81482
$ws.Range("K135").Value = 4444.33338
$ws.Range("M135").Value = -1909.33338
# Row 138
$ws.Range("H138").Value = 2268.26
$ws.Range("I138").Value = 1038.8334
$ws.Range("J138").Value = 2435.9092
$ws.Range("K138").Value = 3116.5002
$ws.Range("L138").Value = 7307.7276
$ws.Range("M138").Value = 2023.4998
$ws.Range("N138").Value = -17587.7276
# Row 141
$ws.Range("H141").Value = 997.5
$ws.Range("I141").Value = 997.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2992.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2187.5
$ws.Range("N141").ClearContents()

# ---------------- Sheet: ARM ----------------
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9097.358
$ws.Range("I32").Value = 6671.531
$ws.Range("K32").Value = 6671.531
$ws.Range("M32").Value = -6384.531
# Row 45
$ws.Range("H45").Value = 1082.625
$ws.Range("I45").Value = 918.5
$ws.Range("J45").Value = 1575
$ws.Range("K45").Value = 918.5
$ws.Range("L45").Value = 1575
$ws.Range("M45").Value = -541.5
$ws.Range("N45").Value = -2329
# Row 57
$ws.Range("H57").Value = 2900
$ws.Range("I57").Value = 2900
$ws.Range("K57").Value = 2900
$ws.Range("M57").Value = -2416
# Row 61
$ws.Range("H61").Value = 71429970
$ws.Range("J61").Value = 1905.2
$ws.Range("L61").Value = 1905.2
$ws.Range("N61").Value = -2329.2
# Row 74
$ws.Range("H74").Value = 3286.6667
$ws.Range("I74").Value = 2546.6667
$ws.Range("K74").Value = 2546.6667
$ws.Range("M74").Value = -1672.6667
# Row 77
$ws.Range("H77").Value = 3286.6667
$ws.Range("I77").Value = 2546.6667
$ws.Range("K77").Value = 12733.3335
$ws.Range("M77").Value = -8365.333500000001
# Row 122
$ws.Range("H122").Value = 3932.0625
$ws.Range("I122").Value = 3784.5386
$ws.Range("J122").Value = 4571.3335
$ws.Range("K122").Value = 11353.6158
$ws.Range("L122").Value = 13714.0005
$ws.Range("M122").Value = -8903.6158
$ws.Range("N122").Value = -18614.0005
# Row 132
$ws.Range("H132").Value = 3244.658
$ws.Range("I132").Value = 2617.0386
$ws.Range("J132").Value = 4604.5
$ws.Range("K132").Value = 7851.1158
$ws.Range("L132").Value = 13813.5
$ws.Range("M132").Value = -5321.1158
$ws.Range("N132").Value = -18873.5
# Row 136
$ws.Range("H136").Value = 71429970
$ws.Range("J136").Value = 1905.2
$ws.Range("L136").Value = 5715.6
$ws.Range("N136").Value = -10815.6

# ---------------- Sheet: BSM ----------------
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 6657.8887
$ws.Range("I134").Value = 1115.1875
$ws.Range("K134").Value = 3345.5625
$ws.Range("M134").Value = -810.5625

# ---------------- Sheet: CRP ----------------
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 327.7857
$ws.Range("I7").Value = 110.77778
$ws.Range("J7").Value = 718.4
$ws.Range("K7").Value = 110.77778
$ws.Range("L7").Value = 718.4
$ws.Range("M7").Value = 2.222219999999993
$ws.Range("N7").Value = -944.4
# Row 58
$ws.Range("H58").Value = 4970.5557
$ws.Range("I58").Value = 1392.2778
$ws.Range("K58").Value = 1392.2778
$ws.Range("M58").Value = -1189.2778
# Row 132
$ws.Range("H132").Value = 2856.5
$ws.Range("I132").Value = 2552.6667
$ws.Range("J132").Value = 3312.25
$ws.Range("K132").Value = 7658.000100000001
$ws.Range("L132").Value = 9936.75
$ws.Range("M132").Value = -5128.000100000001
$ws.Range("N132").Value = -14996.75
# Row 136
$ws.Range("H136").Value = 4970.5557
$ws.Range("I136").Value = 1392.2778
$ws.Range("K136").Value = 4176.8334
$ws.Range("M136").Value = -1626.8334

# ---------------- Sheet: CUL ----------------
$ws = $wb.Worksheets.Item("CUL")
# Row 29
$ws.Range("H29").Value = 770
$ws.Range("I29").Value = 100
$ws.Range("J29").Value = 993.3333
$ws.Range("K29").Value = 300
$ws.Range("L29").Value = 2979.9999
$ws.Range("M29").Value = -23
$ws.Range("N29").Value = -3533.9999
# Row 35
$ws.Range("H35").Value = 4335
$ws.Range("J35").Value = 4335
$ws.Range("L35").Value = 13005
$ws.Range("N35").Value = -13581
# Row 81
$ws.Range("H81").Value = 3466.5833
$ws.Range("J81").Value = 3900
$ws.Range("L81").Value = 11700
$ws.Range("N81").Value = -13946
# Row 84
$ws.Range("H84").Value = 3466.5833
$ws.Range("J84").Value = 3900
$ws.Range("L84").Value = 35100
$ws.Range("N84").Value = -46332
# Row 92
$ws.Range("H92").Value = 519.6667
$ws.Range("J92").Value = 512.5
$ws.Range("L92").Value = 1537.5
$ws.Range("N92").Value = -4033.5
# Row 131
$ws.Range("H131").Value = 20030900
$ws.Range("J131").Value = 42698.527
$ws.Range("L131").Value = 128095.581
$ws.Range("N131").Value = -138175.581
# Row 132
$ws.Range("H132").Value = 1154.4445
$ws.Range("I132").Value = 907.5
$ws.Range("J132").Value = 1463.125
$ws.Range("K132").Value = 8167.5
$ws.Range("L132").Value = 13168.125
$ws.Range("M132").Value = -5637.5
$ws.Range("N132").Value = -18228.125

# ---------------- Sheet: GSM ----------------
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 32146800
$ws.Range("I70").Value = 25004120
$ws.Range("J70").Value = 50003500
$ws.Range("K70").Value = 25004120
$ws.Range("L70").Value = 50003500
$ws.Range("M70").Value = -25003850
$ws.Range("N70").Value = -50004040
# Row 73
$ws.Range("H73").Value = 32146800
$ws.Range("I73").Value = 25004120
$ws.Range("J73").Value = 50003500
$ws.Range("K73").Value = 25004120
$ws.Range("L73").Value = 50003500
$ws.Range("M73").Value = -25003184
$ws.Range("N73").Value = -50005372
# Row 102
$ws.Range("H102").Value = 2327.9092
$ws.Range("I102").Value = 1767
$ws.Range("J102").Value = 3089.1428
$ws.Range("K102").Value = 1767
$ws.Range("L102").Value = 3089.1428
$ws.Range("M102").Value = -145
$ws.Range("N102").Value = -6333.1428
# Row 122
$ws.Range("H122").Value = 967.7273
$ws.Range("I122").Value = 988.7
$ws.Range("J122").Value = 758
$ws.Range("K122").Value = 2966.1
$ws.Range("L122").Value = 2274
$ws.Range("M122").Value = -516.1000000000004
$ws.Range("N122").Value = -7174
# Row 128
$ws.Range("H128").Value = 37000
$ws.Range("I128").Value = 37000
$ws.Range("K128").Value = 37000
$ws.Range("M128").Value = -32020
# Row 132
$ws.Range("H132").Value = 3322.6365
$ws.Range("I132").Value = 3085.4
$ws.Range("J132").Value = 4064
$ws.Range("K132").Value = 9256.200000000001
$ws.Range("L132").Value = 12192
$ws.Range("M132").Value = -6726.200000000001
$ws.Range("N132").Value = -17252

# ---------------- Sheet: LTW ----------------
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2900.4546
$ws.Range("J132").Value = 2388.25
$ws.Range("L132").Value = 7164.75
$ws.Range("N132").Value = -12224.75
# Row 136
$ws.Range("H136").Value = 2200
$ws.Range("I136").Value = 1666.6666
$ws.Range("K136").Value = 4999.9998
$ws.Range("M136").Value = -2449.9998

# ---------------- Sheet: WVR ----------------
$ws = $wb.Worksheets.Item("WVR")
# Row 25
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
# Row 107
$ws.Range("H107").Value = 473
$ws.Range("I107").Value = 473
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1419
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 501
$ws.Range("N107").ClearContents()
# Row 132
$ws.Range("H132").Value = 2996.32
$ws.Range("I132").Value = 2754
$ws.Range("J132").Value = 3359.8
$ws.Range("K132").Value = 8262
$ws.Range("L132").Value = 10079.4
$ws.Range("M132").Value = -5732
$ws.Range("N132").Value = -15139.4
# Row 136
$ws.Range("H136").Value = 1410.6111
$ws.Range("I136").Value = 1099.4
$ws.Range("J136").Value = 2966.6667
$ws.Range("K136").Value = 3298.2
$ws.Range("L136").Value = 8900.000100000001
$ws.Range("M136").Value = -748.2000000000003
$ws.Range("N136").Value = -14000.0001
